$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$rng = $ws.Range("B28")
$formulaText = '"if($E$7="""";true;false)"'
$fc1 = $rng.FormatConditions.Add(2, 3, $formulaText)
$fc2 = $rng.FormatConditions.AddIconSetCondition()
$fc2.IconSet = 13
$fc1.SetFirstPriority()
$fc2.SetFirstPriority()
